# Applies the odds updates described in the commit diff to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.62
$ws.Range("V2").Value = 1.54

# Row 3
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("Q3").Value = 2.5
$ws.Range("R3").Value = 1.5
$ws.Range("V3").Value = 1.63
$ws.Range("BC3").Value = 126

# Row 4
$ws.Range("I4").Value = 3.1
$ws.Range("J4").Value = 3.6
$ws.Range("V4").Value = 1.47

# Row 5
$ws.Range("I5").Value = 7.5
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.75
$ws.Range("Q5").Value = 1.8
$ws.Range("U5").Value = 1.87
$ws.Range("V5").Value = 1.77
$ws.Range("W5").Value = 6.5
$ws.Range("AH5").Value = 41
$ws.Range("AI5").Value = 23
$ws.Range("AN5").Value = 7
$ws.Range("AR5").Value = 151
